# ------------------------------------------------------------------
# Apply the "imported JSON type test case from Excel" edit:
#   1. Rename Sheet2 -> RegisterTestCaseSheet2
#   2. sheet1 (RegisterTestCaseSheet): tweak selection + F6 value
#   3. sheet2 (RegisterTestCaseSheet2): populate with a mirrored table
#      (A-E copied from sheet1) plus a new "params" (JSON) column F,
#      matching hyperlinks/styles/column widths/selection.
# ------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- workbook: rename the second sheet -----------------------------
$ws2.Name = "RegisterTestCaseSheet2"

# --- sheet1: small tweaks --------------------------------------------
$ws1.Range("F6").Value = 13344445355

# --- sheet2: header row ----------------------------------------------
$ws2.Cells.Item(1, 1).Value = "CaseId"
$ws2.Cells.Item(1, 2).Value = "TestName(接口名）"
$ws2.Cells.Item(1, 3).Value = "Url"
$ws2.Cells.Item(1, 4).Value = "Type"
$ws2.Cells.Item(1, 5).Value = "Desc(用例描述）"
# "params" must become shared-string index 16, and F7's value must land
# right after it (index 17) before F2..F6 get appended - this mirrors
# the order the strings ended up in the source workbook.
$ws2.Cells.Item(1, 6).Value = "params"

# --- sheet2: data rows (A-E mirror sheet1's rows 2-7) -----------------
$ws2.Cells.Item(2, 1).Value = 1
$ws2.Cells.Item(2, 2).Value = "Register"
$ws2.Cells.Item(2, 3).Value = "http://47.107.166.132:8080/futureloan/mvc/api/member/register"
$ws2.Cells.Item(2, 4).Value = "PSOT"
$ws2.Cells.Item(2, 5).Value = "有手机号，无密码"

$ws2.Cells.Item(3, 1).Value = 2
$ws2.Cells.Item(3, 2).Value = "Register"
$ws2.Cells.Item(3, 3).Value = "http://47.107.166.132:8080/futureloan/mvc/api/member/register"
$ws2.Cells.Item(3, 4).Value = "PSOT"
$ws2.Cells.Item(3, 5).Value = "无手机号，有密码"

$ws2.Cells.Item(4, 1).Value = 3
$ws2.Cells.Item(4, 2).Value = "Register"
$ws2.Cells.Item(4, 3).Value = "http://47.107.166.132:8080/futureloan/mvc/api/member/register"
$ws2.Cells.Item(4, 4).Value = "PSOT"
$ws2.Cells.Item(4, 5).Value = "不合格手机号"

$ws2.Cells.Item(5, 1).Value = 4
$ws2.Cells.Item(5, 2).Value = "Register"
$ws2.Cells.Item(5, 3).Value = "http://47.107.166.132:8080/futureloan/mvc/api/member/register"
$ws2.Cells.Item(5, 4).Value = "PSOT"
$ws2.Cells.Item(5, 5).Value = "合格手机号+不合格密码"

$ws2.Cells.Item(6, 1).Value = 5
$ws2.Cells.Item(6, 2).Value = "Register"
$ws2.Cells.Item(6, 3).Value = "http://47.107.166.132:8080/futureloan/mvc/api/member/register"
$ws2.Cells.Item(6, 4).Value = "PSOT"
$ws2.Cells.Item(6, 5).Value = "合格手机号+合格密码"

$ws2.Cells.Item(7, 1).Value = 6
$ws2.Cells.Item(7, 2).Value = "Register"
$ws2.Cells.Item(7, 3).Value = "http://47.107.166.132:8080/futureloan/mvc/api/member/register"
$ws2.Cells.Item(7, 4).Value = "PSOT"
$ws2.Cells.Item(7, 5).Value = "重复上述5"

# --- sheet2: F column (new JSON "params") ------------------------------
# Write row 7 first so it becomes shared-string #17 (right after "params"
# at #16), then rows 2-6 in order (#18..#22) - reproduces the exact
# shared-string table ordering of the authored workbook.
$ws2.Cells.Item(7, 6).Value = '{"mobilephone":"13344445555","pwd":"123456"}'
$ws2.Cells.Item(2, 6).Value = '{"mobilephone":"13344445555","pwd":""}'
$ws2.Cells.Item(3, 6).Value = '{"mobilephone":"","pwd":"123456"}'
$ws2.Cells.Item(4, 6).Value = '{"mobilephone":"1334444","pwd":"123456"}'
$ws2.Cells.Item(5, 6).Value = '{"mobilephone":"13344445555","pwd":"12345"}'
$ws2.Cells.Item(6, 6).Value = '{"mobilephone":"13344445545","pwd":"123456"}'

# --- sheet2: alignment (center), matching sheet1's cellXfs s="1" -------
# (B1 intentionally stays un-styled, mirroring sheet1's quirk.)
$ws2.Range("A1").HorizontalAlignment = -4108
$ws2.Range("C1:G7").HorizontalAlignment = -4108
$ws2.Range("A2:B7").HorizontalAlignment = -4108

# --- sheet2: hyperlinks on the Url column ------------------------------
$ws2.Hyperlinks.Add($ws2.Range("C3:C7"), "http://47.107.166.132:8080/futureloan/mvc/api/member/register", "", "", "http://47.107.166.132:8080/futureloan/mvc/api/member/register")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "http://47.107.166.132:8080/futureloan/mvc/api/member/register", "", "", "http://47.107.166.132:8080/futureloan/mvc/api/member/register")
$ws2.Range("C2:C7").HorizontalAlignment = -4108

# --- sheet2: column widths (converted to this engine's MDW=6 rounding) -
$ws2.Columns.Item(1).ColumnWidth = 18.666666666666668
$ws2.Columns.Item(2).ColumnWidth = 22.666666666666668
$ws2.Columns.Item(3).ColumnWidth = 62.333333333333336
$ws2.Columns.Item(4).ColumnWidth = 18.833333333333332
$ws2.Columns.Item(5).ColumnWidth = 23.166666666666668
$ws2.Columns.Item(6).ColumnWidth = 48.5
$ws2.Columns.Item(7).ColumnWidth = 17.333333333333332

# --- selections ---------------------------------------------------------
$ws2.Range("C27").Select()
$ws1.Range("E12").Select()
